$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update nombre_aides (column C) and montant_total (column D) for the
# rows affected by the 2020-07-16 data refresh.
$ws.Range("C2").Value = 34801
$ws.Range("D2").Value = 50405229
$ws.Range("C3").Value = 85291
$ws.Range("D3").Value = 125194582
$ws.Range("C4").Value = 29249
$ws.Range("D4").Value = 43360493
$ws.Range("C5").Value = 8047
$ws.Range("D5").Value = 11965946
$ws.Range("C6").Value = 1719
$ws.Range("D6").Value = 2558542
$ws.Range("C11").Value = 38214
$ws.Range("D11").Value = 52000789
$ws.Range("C12").Value = 8956
$ws.Range("D12").Value = 12968630
$ws.Range("C13").Value = 24565
$ws.Range("D13").Value = 36056867
$ws.Range("C14").Value = 7831
$ws.Range("D14").Value = 11633050
$ws.Range("C15").Value = 1985
$ws.Range("D15").Value = 2952976
$ws.Range("C19").Value = 9490
$ws.Range("D19").Value = 12617260
$ws.Range("C20").Value = 12513
$ws.Range("D20").Value = 18081769
$ws.Range("C21").Value = 29986
$ws.Range("D21").Value = 44061596
$ws.Range("C22").Value = 9739
$ws.Range("D22").Value = 14485922
$ws.Range("C23").Value = 2461
$ws.Range("D23").Value = 3661763
$ws.Range("C26").Value = 10928
$ws.Range("D26").Value = 14671795
$ws.Range("C27").Value = 7122
$ws.Range("D27").Value = 10325993
$ws.Range("C28").Value = 21235
$ws.Range("D28").Value = 31200515
$ws.Range("C29").Value = 7386
$ws.Range("D29").Value = 10994651
$ws.Range("C30").Value = 1816
$ws.Range("D30").Value = 2712486
$ws.Range("C33").Value = 7732
$ws.Range("D33").Value = 10253101
$ws.Range("C34").Value = 2776
$ws.Range("D34").Value = 4004910
$ws.Range("C35").Value = 6989
$ws.Range("D35").Value = 10213847
$ws.Range("C36").Value = 2833
$ws.Range("D36").Value = 4194523
$ws.Range("C40").Value = 2152
$ws.Range("D40").Value = 2902328
$ws.Range("C41").Value = 16029
$ws.Range("D41").Value = 23199910
$ws.Range("C42").Value = 48185
$ws.Range("D42").Value = 70702622
$ws.Range("C43").Value = 18050
$ws.Range("D43").Value = 26818655
$ws.Range("C44").Value = 5237
$ws.Range("D44").Value = 7806108
$ws.Range("C45").Value = 1033
$ws.Range("D45").Value = 1541684
$ws.Range("C49").Value = 15479
$ws.Range("D49").Value = 20681898
$ws.Range("C50").Value = 1675
$ws.Range("D50").Value = 2431133
$ws.Range("C51").Value = 5947
$ws.Range("D51").Value = 8756326
$ws.Range("C52").Value = 2083
$ws.Range("D52").Value = 3111370
$ws.Range("C53").Value = 678
$ws.Range("D53").Value = 1012305
$ws.Range("C56").Value = 5408
$ws.Range("D56").Value = 7473915
$ws.Range("C57").Value = 679
$ws.Range("D57").Value = 994840
$ws.Range("C58").Value = 1722
$ws.Range("D58").Value = 2551414
$ws.Range("C59").Value = 692
$ws.Range("D59").Value = 1032219
$ws.Range("C62").Value = 12
$ws.Range("D62").Value = 18000
$ws.Range("C63").Value = 1001
$ws.Range("D63").Value = 1419265
$ws.Range("C64").Value = 14282
$ws.Range("D64").Value = 20650106
$ws.Range("C65").Value = 42321
$ws.Range("D65").Value = 61982652
$ws.Range("C66").Value = 14936
$ws.Range("D66").Value = 22210747
$ws.Range("C67").Value = 4291
$ws.Range("D67").Value = 6393293
$ws.Range("C68").Value = 825
$ws.Range("D68").Value = 1228273
$ws.Range("C71").Value = 14160
$ws.Range("D71").Value = 18751339
$ws.Range("C72").Value = 45483
$ws.Range("D72").Value = 66229049
$ws.Range("C73").Value = 132096
$ws.Range("D73").Value = 194758050
$ws.Range("C74").Value = 58127
$ws.Range("D74").Value = 86651715
$ws.Range("C75").Value = 18436
$ws.Range("D75").Value = 27552074
$ws.Range("C76").Value = 4095
$ws.Range("D76").Value = 6119120
$ws.Range("C83").Value = 44956
$ws.Range("D83").Value = 61519156
$ws.Range("C84").Value = 4077
$ws.Range("D84").Value = 5913759
$ws.Range("C85").Value = 10570
$ws.Range("D85").Value = 15537460
$ws.Range("C86").Value = 3622
$ws.Range("D86").Value = 5399081
$ws.Range("C87").Value = 1261
$ws.Range("D87").Value = 1886091
$ws.Range("C88").Value = 257
$ws.Range("D88").Value = 383612
$ws.Range("C91").Value = 4742
$ws.Range("D91").Value = 6398037
$ws.Range("C92").Value = 1378
$ws.Range("D92").Value = 1991490
$ws.Range("C93").Value = 4571
$ws.Range("D93").Value = 6735995
$ws.Range("C94").Value = 1775
$ws.Range("D94").Value = 2647116
$ws.Range("C95").Value = 622
$ws.Range("D95").Value = 932141
$ws.Range("C99").Value = 3028
$ws.Range("D99").Value = 4020458
$ws.Range("C101").Value = 263
$ws.Range("D101").Value = 392765
$ws.Range("C102").Value = 92
$ws.Range("D102").Value = 138000
$ws.Range("C105").Value = 10022
$ws.Range("D105").Value = 14563734
$ws.Range("C106").Value = 27718
$ws.Range("D106").Value = 40753368
$ws.Range("C107").Value = 9296
$ws.Range("D107").Value = 13827476
$ws.Range("C108").Value = 2521
$ws.Range("D108").Value = 3758910
$ws.Range("C112").Value = 9137
$ws.Range("D112").Value = 12111675
$ws.Range("C113").Value = 28172
$ws.Range("D113").Value = 40677284
$ws.Range("C114").Value = 62450
$ws.Range("D114").Value = 91483828
$ws.Range("C115").Value = 20238
$ws.Range("D115").Value = 30097470
$ws.Range("C116").Value = 5657
$ws.Range("D116").Value = 8432861
$ws.Range("C117").Value = 1009
$ws.Range("D117").Value = 1509493
$ws.Range("C121").Value = 24003
$ws.Range("D121").Value = 32162851
$ws.Range("C122").Value = 32842
$ws.Range("D122").Value = 47458653
$ws.Range("C123").Value = 71611
$ws.Range("D123").Value = 104836180
$ws.Range("C124").Value = 22413
$ws.Range("D124").Value = 33280082
$ws.Range("C125").Value = 5919
$ws.Range("D125").Value = 8803044
$ws.Range("C126").Value = 1073
$ws.Range("D126").Value = 1599769
$ws.Range("C130").Value = 29097
$ws.Range("D130").Value = 38774729
$ws.Range("C131").Value = 12294
$ws.Range("D131").Value = 17806921
$ws.Range("C132").Value = 30641
$ws.Range("D132").Value = 45042386
$ws.Range("C133").Value = 10917
$ws.Range("D133").Value = 16224584
$ws.Range("C134").Value = 2757
$ws.Range("D134").Value = 4111791
$ws.Range("C135").Value = 438
$ws.Range("D135").Value = 650990
$ws.Range("C138").Value = 10120
$ws.Range("D138").Value = 13562697
$ws.Range("C139").Value = 32061
$ws.Range("D139").Value = 46334898
$ws.Range("C140").Value = 75817
$ws.Range("D140").Value = 111159024
$ws.Range("C141").Value = 22883
$ws.Range("D141").Value = 34029298
$ws.Range("C142").Value = 5893
$ws.Range("D142").Value = 8796848
$ws.Range("C143").Value = 1269
$ws.Range("D143").Value = 1891982
$ws.Range("C146").Value = 26992
$ws.Range("D146").Value = 36586701
